$d = $word.ActiveDocument

function Replace-ParagraphInner {
    param($Paragraph, $InnerXml)
    $rng = $Paragraph.Range
    # Clear everything except the trailing paragraph mark, then insert the
    # replacement run content (with proofErr markers etc.) in its place.
    $inner = $d.Range($rng.Start, $rng.End - 1)
    $inner.Text = ""
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        $InnerXml +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $inner.InsertXML($xml)
}

$rPr24 = '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

# --- Paragraph "mkdir ... make directory" -> wrap "mkdir" in spellcheck markers ---
$pMkdir = $d.Paragraphs.Item(3)
$innerMkdir =
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr24 + '<w:t>mkdir</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="00122151">' + $rPr24 + '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> make directory</w:t></w:r>'
Replace-ParagraphInner $pMkdir $innerMkdir

# --- Paragraph "git init ... initialize an empty git repository..." -> wrap "init" ---
$pGitInit = $d.Paragraphs.Item(6)
$innerGitInit =
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">git </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPr24 + '<w:t>init</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="00122151">' + $rPr24 + '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> initialize an empty git repository in particular folder</w:t></w:r>'
Replace-ParagraphInner $pGitInit $innerGitInit

# --- Paragraph "git status ... where .git folder is there" -> wrap "where .git" ---
$pGitStatus = $d.Paragraphs.Item(8)
$innerGitStatus =
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">git status </w:t></w:r>' +
    '<w:r w:rsidRPr="00EF2144">' + $rPr24 + '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> will show untracked files only in that directory </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPr24 + '<w:t>where .git</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> folder is there</w:t></w:r>'
Replace-ParagraphInner $pGitStatus $innerGitStatus

# --- Paragraph "git add . ... all files will be added" -> wrap "add ." ---
$pGitAdd = $d.Paragraphs.Item(9)
$innerGitAdd =
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve">git </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r>' + $rPr24 + '<w:t>add .</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r w:rsidRPr="00EF2144">' + $rPr24 + '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
    '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> all files will be added</w:t></w:r>'
Replace-ParagraphInner $pGitAdd $innerGitAdd

# --- Append the three new paragraphs at the end of the document ---
$endRng = $d.Content
$endRng.Collapse(0)

$newParasXml =
    '<w:p><w:pPr>' + $rPr24 + '</w:pPr>' +
        '<w:r>' + $rPr24 + '<w:t xml:space="preserve">git log </w:t></w:r>' +
        '<w:r>' + $rPr24 + '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
        '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> entire history of the project commit' + [char]0x2019 + 's</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/>' +
        '<w:color w:val="C0A000"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
        '</w:pPr>' +
        '<w:r>' + $rPr24 + '<w:t xml:space="preserve">git reset </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/>' +
            '<w:color w:val="C0A000"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
            '<w:t>2b82c969868f8cebeb4fdbb879e5e3220771629f</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/>' +
            '<w:color w:val="C0A000"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
            '<w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Lucida Console" w:hAnsi="Lucida Console" w:cs="Lucida Console"/>' +
            '<w:color w:val="C0A000"/><w:kern w:val="0"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr>' +
            '<w:sym w:font="Wingdings" w:char="F0E0"/></w:r>' +
        '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> command to go to that particular change (using hash code), changes above </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r>' + $rPr24 + '<w:t>this  change</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r>' + $rPr24 + '<w:t xml:space="preserve"> will be removed</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr>' + $rPr24 + '</w:pPr></w:p>'

$xmlPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    $newParasXml +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$endRng.InsertXML($xmlPkg)

Write-Output "edit complete"
